$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data: Stoyan, 23 (Stoyan already exists in shared strings from row 3)
$ws.Range("A6").Value = "Stoyan"
$ws.Range("B6").Value = 23

# Match formatting of existing rows (A: left aligned, B: right aligned)
$ws.Range("A6").HorizontalAlignment = -4131  # xlLeft
$ws.Range("B6").HorizontalAlignment = -4152  # xlRight

# Update the active selection to A7, as shown in the diff
$ws.Range("A7").Select()
